$d = $word.ActiveDocument

$replacements = @(
    @("36×19=", "37×47="),
    @("32×19=", "24×79="),
    @("47×58=", "21×53="),
    @("77×41=", "92×23="),
    @("37×33=", "88×88="),
    @("12×98=", "92×17="),
    @("40×35=", "70×82="),
    @("56×72=", "73×57="),
    @("75×40=", "35×20="),
    @("46×82=", "77×57="),
    @("48×61=", "14×12="),
    @("52×56=", "18×70="),
    @("69×67=", "31×78="),
    @("48×15=", "61×47="),
    @("66×75=", "83×74="),
    @("59×82=", "87×18="),
    @("79×87=", "42×37="),
    @("93×82=", "96×38="),
    @("90×48=", "24×13="),
    @("91×75=", "23×68="),
    @("36×39=", "52×69="),
    @("84×75=", "73×30="),
    @("60×46=", "54×51="),
    @("71×91=", "77×47="),
    @("18×46=", "38×54=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
